$d = $word.ActiveDocument

$replacements = @(
    @("154×8=1232", "735×3=2205"),
    @("386×7=2702", "205×9=1845"),
    @("254×7=1778", "585×2=1170"),
    @("601×9=5409", "403×2=806"),
    @("825×5=4125", "867×9=7803"),
    @("169×4=676", "229×2=458"),
    @("920×6=5520", "124×8=992"),
    @("874×2=1748", "589×9=5301"),
    @("987×2=1974", "744×2=1488"),
    @("200×6=1200", "415×7=2905"),
    @("636×4=2544", "167×4=668"),
    @("266×7=1862", "655×3=1965"),
    @("323×3=969", "445×8=3560"),
    @("400×2=800", "468×4=1872"),
    @("119×4=476", "223×7=1561"),
    @("169×2=338", "603×2=1206"),
    @("312×4=1248", "809×9=7281"),
    @("491×6=2946", "808×8=6464"),
    @("598×4=2392", "251×3=753"),
    @("338×6=2028", "233×5=1165"),
    @("878×4=3512", "325×8=2600"),
    @("998×9=8982", "210×2=420"),
    @("621×5=3105", "526×8=4208"),
    @("588×4=2352", "330×9=2970"),
    @("941×6=5646", "284×2=568")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

Write-Host "Done"
